$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.039.96'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.29%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.712.68'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.11%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9979'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.52%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '318.36'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.23%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9996'
$ws.Range("D6").Style = "Normal"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4040'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.77%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4088'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.59%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.486'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.85%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '54.08'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.83%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.001'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.12%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08848'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.07%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '26.45'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +6.39%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.533'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.44%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.152'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.09%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001364'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.63%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.721.07'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.51%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '97.23'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.28%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.07172'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.44%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '21.25'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.38%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.292'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.82%  '

$ws.Range("E22").Value = '  -0.47%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.43'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.89%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '24.953.32'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.01%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.332'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.07%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.922'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -7.02%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '23.36'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.70%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.271'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +19.98%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '167.42'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.19%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '146.97'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +4.65%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.410'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -9.74%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.920.53'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.92%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.230'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +13.45%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08907'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.47%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.03240'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +7.04%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '7.293'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -7.65%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.033'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -5.38%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2869'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.99%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.8515'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.77%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '10.92'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.85%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.09368'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.61%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '14.26'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.31%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.473'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.13%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '17.48'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.44%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.735'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.27%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.7471'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.15%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.248'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.56%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.403'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.80%  '

$ws.Range("E49").Value = '  -0.19%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '142.51'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.02%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.08397'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.70%  '

